$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.551825333333333
$ws.Range("H2").Value = 13.655476
$ws.Range("I2").Value = 0.3901863008207799
$ws.Range("J2").Value = 0.3901863008207799
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.131233999999999
$ws.Range("N2").Value = 24.393702
$ws.Range("O2").Value = 0.02090995573015822
$ws.Range("P2").Value = 0.02090995573015823
$ws.Range("Q2").Value = 37.01195691246133
$ws.Range("R2").Value = 333.107612212152
$ws.Range("S2").Value = 0.008158778276676708
$ws.Range("T2").Value = 0.00815877827667671
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.551825333333333
$ws.Range("H3").Value = 13.655476
$ws.Range("I3").Value = 0.3901863008207799
$ws.Range("J3").Value = 0.3901863008207799
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 243.3763986666667
$ws.Range("N3").Value = 730.1291960000001
$ws.Range("O3").Value = 0.625857000534647
$ws.Range("P3").Value = 0.6258570005346471
$ws.Range("Q3").Value = 1107.806856986366
$ws.Range("R3").Value = 9970.261712877298
$ws.Range("S3").Value = 0.2442008278814028
$ws.Range("T3").Value = 0.2442008278814028
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.551825333333333
$ws.Range("H4").Value = 13.655476
$ws.Range("I4").Value = 0.3901863008207799
$ws.Range("J4").Value = 0.3901863008207799
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 103.9426383333333
$ws.Range("N4").Value = 311.827915
$ws.Range("O4").Value = 0.2672947262403034
$ws.Range("P4").Value = 0.2672947262403035
$ws.Range("Q4").Value = 473.1287343791711
$ws.Range("R4").Value = 4258.15860941254
$ws.Range("S4").Value = 0.104294740460607
$ws.Range("T4").Value = 0.1042947404606071
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.551825333333333
$ws.Range("H5").Value = 13.655476
$ws.Range("I5").Value = 0.3901863008207799
$ws.Range("J5").Value = 0.3901863008207799
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 33.41874933333333
$ws.Range("N5").Value = 100.256248
$ws.Range("O5").Value = 0.08593831749489127
$ws.Range("P5").Value = 0.08593831749489128
$ws.Range("Q5").Value = 152.1163098237831
$ws.Range("R5").Value = 1369.046788414048
$ws.Range("S5").Value = 0.03353195420209334
$ws.Range("T5").Value = 0.03353195420209334
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3.622039333333333
$ws.Range("H6").Value = 10.866118
$ws.Range("I6").Value = 0.3104842618962599
$ws.Range("J6").Value = 0.3104842618962599
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.131233999999999
$ws.Range("N6").Value = 24.393702
$ws.Range("O6").Value = 0.02090995573015822
$ws.Range("P6").Value = 0.02090995573015823
$ws.Range("Q6").Value = 29.45164937653733
$ws.Range("R6").Value = 265.064844388836
$ws.Range("S6").Value = 0.006492212171161645
$ws.Range("T6").Value = 0.006492212171161647
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3.622039333333333
$ws.Range("H7").Value = 10.866118
$ws.Range("I7").Value = 0.3104842618962599
$ws.Range("J7").Value = 0.3104842618962599
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 243.3763986666667
$ws.Range("N7").Value = 730.1291960000001
$ws.Range("O7").Value = 0.625857000534647
$ws.Range("P7").Value = 0.6258570005346471
$ws.Range("Q7").Value = 881.518888775681
$ws.Range("R7").Value = 7933.66999898113
$ws.Range("S7").Value = 0.194318748863607
$ws.Range("T7").Value = 0.194318748863607
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.622039333333333
$ws.Range("H8").Value = 10.866118
$ws.Range("I8").Value = 0.3104842618962599
$ws.Range("J8").Value = 0.3104842618962599
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 103.9426383333333
$ws.Range("N8").Value = 311.827915
$ws.Range("O8").Value = 0.2672947262403034
$ws.Range("P8").Value = 0.2672947262403035
$ws.Range("Q8").Value = 376.4843244537744
$ws.Range("R8").Value = 3388.35892008397
$ws.Range("S8").Value = 0.08299080578548346
$ws.Range("T8").Value = 0.08299080578548347
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.622039333333333
$ws.Range("H9").Value = 10.866118
$ws.Range("I9").Value = 0.3104842618962599
$ws.Range("J9").Value = 0.3104842618962599
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 33.41874933333333
$ws.Range("N9").Value = 100.256248
$ws.Range("O9").Value = 0.08593831749489127
$ws.Range("P9").Value = 0.08593831749489128
$ws.Range("Q9").Value = 121.0440245561404
$ws.Range("R9").Value = 1089.396221005264
$ws.Range("S9").Value = 0.02668249507600775
$ws.Range("T9").Value = 0.02668249507600776
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.291101666666667
$ws.Range("H10").Value = 3.873305
$ws.Range("I10").Value = 0.1106743221474397
$ws.Range("J10").Value = 0.1106743221474397
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.131233999999999
$ws.Range("N10").Value = 24.393702
$ws.Range("O10").Value = 0.02090995573015822
$ws.Range("P10").Value = 0.02090995573015823
$ws.Range("Q10").Value = 10.49824976945667
$ws.Range("R10").Value = 94.48424792511
$ws.Range("S10").Value = 0.002314195176568234
$ws.Range("T10").Value = 0.002314195176568234
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.291101666666667
$ws.Range("H11").Value = 3.873305
$ws.Range("I11").Value = 0.1106743221474397
$ws.Range("J11").Value = 0.1106743221474397
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 243.3763986666667
$ws.Range("N11").Value = 730.1291960000001
$ws.Range("O11").Value = 0.625857000534647
$ws.Range("P11").Value = 0.6258570005346471
$ws.Range("Q11").Value = 314.2236739458645
$ws.Range("R11").Value = 2828.01306551278
$ws.Range("S11").Value = 0.06926629929540186
$ws.Range("T11").Value = 0.06926629929540186
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.291101666666667
$ws.Range("H12").Value = 3.873305
$ws.Range("I12").Value = 0.1106743221474397
$ws.Range("J12").Value = 0.1106743221474397
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 103.9426383333333
$ws.Range("N12").Value = 311.827915
$ws.Range("O12").Value = 0.2672947262403034
$ws.Range("P12").Value = 0.2672947262403035
$ws.Range("Q12").Value = 134.2005135898972
$ws.Range("R12").Value = 1207.804622309075
$ws.Range("S12").Value = 0.02958266264023104
$ws.Range("T12").Value = 0.02958266264023104
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.291101666666667
$ws.Range("H13").Value = 3.873305
$ws.Range("I13").Value = 0.1106743221474397
$ws.Range("J13").Value = 0.1106743221474397
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 33.41874933333333
$ws.Range("N13").Value = 100.256248
$ws.Range("O13").Value = 0.08593831749489127
$ws.Range("P13").Value = 0.08593831749489128
$ws.Range("Q13").Value = 43.14700296218222
$ws.Range("R13").Value = 388.32302665964
$ws.Range("S13").Value = 0.009511165035238548
$ws.Range("T13").Value = 0.009511165035238548
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 2.200808
$ws.Range("H14").Value = 6.602424000000001
$ws.Range("I14").Value = 0.1886551151355205
$ws.Range("J14").Value = 0.1886551151355205
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 8.131233999999999
$ws.Range("N14").Value = 24.393702
$ws.Range("O14").Value = 0.02090995573015822
$ws.Range("P14").Value = 0.02090995573015823
$ws.Range("Q14").Value = 17.895284837072
$ws.Range("R14").Value = 161.057563533648
$ws.Range("S14").Value = 0.003944770105751637
$ws.Range("T14").Value = 0.003944770105751638
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 2.200808
$ws.Range("H15").Value = 6.602424000000001
$ws.Range("I15").Value = 0.1886551151355205
$ws.Range("J15").Value = 0.1886551151355205
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 243.3763986666667
$ws.Range("N15").Value = 730.1291960000001
$ws.Range("O15").Value = 0.625857000534647
$ws.Range("P15").Value = 0.6258570005346471
$ws.Range("Q15").Value = 535.6247251967894
$ws.Range("R15").Value = 4820.622526771105
$ws.Range("S15").Value = 0.1180711244942354
$ws.Range("T15").Value = 0.1180711244942354
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 2.200808
$ws.Range("H16").Value = 6.602424000000001
$ws.Range("I16").Value = 0.1886551151355205
$ws.Range("J16").Value = 0.1886551151355205
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 103.9426383333333
$ws.Range("N16").Value = 311.827915
$ws.Range("O16").Value = 0.2672947262403034
$ws.Range("P16").Value = 0.2672947262403035
$ws.Range("Q16").Value = 228.7577899851067
$ws.Range("R16").Value = 2058.820109865961
$ws.Range("S16").Value = 0.05042651735398188
$ws.Range("T16").Value = 0.05042651735398189
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 2.200808
$ws.Range("H17").Value = 6.602424000000001
$ws.Range("I17").Value = 0.1886551151355205
$ws.Range("J17").Value = 0.1886551151355205
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 33.41874933333333
$ws.Range("N17").Value = 100.256248
$ws.Range("O17").Value = 0.08593831749489127
$ws.Range("P17").Value = 0.08593831749489128
$ws.Range("Q17").Value = 73.54825088279468
$ws.Range("R17").Value = 661.9342579451521
$ws.Range("S17").Value = 0.01621270318155163
$ws.Range("T17").Value = 0.01621270318155163
